$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-11 16:18:53'
$ws.Range('H2').Value = "'77%"
$ws.Range('I2').Value = '1.8 mm'
$ws.Range('K2').Value = '5.2 MJ/m2'
$ws.Range('E3').Value = '2026-02-11 16:18:56'
$ws.Range('K3').Value = '8.6 MJ/m2'
$ws.Range('O3').Value = '-0.1 °C'
$ws.Range('E4').Value = '2026-02-11 16:18:59'
$ws.Range('H4').Value = "'58%"
$ws.Range('J4').Value = '1003.1 hPa'
$ws.Range('K4').Value = '10.4 MJ/m2'
$ws.Range('O4').Value = '15.7 °C'
$ws.Range('E5').Value = '2026-02-11 16:19:01'
$ws.Range('H5').Value = "'78%"
$ws.Range('K5').Value = '6.6 MJ/m2'
$ws.Range('E6').Value = '2026-02-11 16:19:04'
$ws.Range('J6').Value = '1003.6 hPa'
$ws.Range('K6').Value = '10.2 MJ/m2'
$ws.Range('O6').Value = '13.1 °C'
$ws.Range('E7').Value = '2026-02-11 16:19:07'
$ws.Range('H7').Value = "'41%"
$ws.Range('J7').Value = '1003.9 hPa'
$ws.Range('K7').Value = '9.3 MJ/m2'
$ws.Range('E8').Value = '2026-02-11 16:19:10'
$ws.Range('K8').Value = '9.9 MJ/m2'
$ws.Range('E9').Value = '2026-02-11 16:19:12'
$ws.Range('K9').Value = '10.3 MJ/m2'
$ws.Range('O9').Value = '12.1 °C'
$ws.Range('E10').Value = '2026-02-11 16:19:15'
$ws.Range('K10').Value = '7.0 MJ/m2'
$ws.Range('O10').Value = '13.4 °C'
$ws.Range('E11').Value = '2026-02-11 16:19:18'
$ws.Range('H11').Value = "'81%"
$ws.Range('O11').Value = '7.5 °C'
$ws.Range('E12').Value = '2026-02-11 16:19:21'
$ws.Range('H12').Value = "'100%"
$ws.Range('O12').Value = '11.4 °C'
$ws.Range('E13').Value = '2026-02-11 16:19:23'
$ws.Range('H13').Value = "'78%"
$ws.Range('J13').Value = '1005.6 hPa'
$ws.Range('K13').Value = '10.4 MJ/m2'
$ws.Range('O13').Value = '7.5 °C'
$ws.Range('E14').Value = '2026-02-11 16:19:26'
$ws.Range('K14').Value = '10.3 MJ/m2'
$ws.Range('E15').Value = '2026-02-11 16:19:28'
$ws.Range('O15').Value = '12.0 °C'
$ws.Range('E16').Value = '2026-02-11 16:19:31'
$ws.Range('I16').Value = '3.9 mm'
$ws.Range('K16').Value = '7.2 MJ/m2'
$ws.Range('E17').Value = '2026-02-11 16:19:33'
$ws.Range('K17').Value = '10.5 MJ/m2'
$ws.Range('E18').Value = '2026-02-11 16:19:36'
$ws.Range('H18').Value = "'73%"
$ws.Range('J18').Value = '1003.6 hPa'
$ws.Range('K18').Value = '10.7 MJ/m2'
$ws.Range('O18').Value = '13.7 °C'
$ws.Range('E19').Value = '2026-02-11 16:19:39'
$ws.Range('H19').Value = "'79%"
$ws.Range('K19').Value = '10.2 MJ/m2'
$ws.Range('L19').Value = '45.7 km/h - 228º 15:55 TU'
$ws.Range('E20').Value = '2026-02-11 16:19:42'
$ws.Range('K20').Value = '9.1 MJ/m2'
$ws.Range('E21').Value = '2026-02-11 16:19:44'
$ws.Range('J21').Value = '1006.2 hPa'
$ws.Range('K21').Value = '9.1 MJ/m2'
$ws.Range('O21').Value = '8.3 °C'
$ws.Range('E22').Value = '2026-02-11 16:19:46'
$ws.Range('G22').Value = '118 cm'
$ws.Range('H22').Value = "'90%"
$ws.Range('K22').Value = '9.2 MJ/m2'
$ws.Range('E23').Value = '2026-02-11 16:19:49'
$ws.Range('H23').Value = "'68%"
$ws.Range('I23').Value = '3.4 mm'
$ws.Range('K23').Value = '7.4 MJ/m2'
$ws.Range('O23').Value = '-1.0 °C'
$ws.Range('E24').Value = '2026-02-11 16:19:52'
$ws.Range('H24').Value = "'69%"
$ws.Range('I24').Value = '4.1 mm'
$ws.Range('O24').Value = '13.7 °C'
$ws.Range('E25').Value = '2026-02-11 16:19:54'
$ws.Range('H25').Value = "'60%"
$ws.Range('K25').Value = '11.1 MJ/m2'
$ws.Range('E26').Value = '2026-02-11 16:19:57'
$ws.Range('J26').Value = '1003.6 hPa'
$ws.Range('K26').Value = '10.5 MJ/m2'
$ws.Range('O26').Value = '6.9 °C'
$ws.Range('E27').Value = '2026-02-11 16:20:00'
$ws.Range('K27').Value = '8.7 MJ/m2'
$ws.Range('E28').Value = '2026-02-11 16:20:03'
$ws.Range('J28').Value = '1003.9 hPa'
$ws.Range('K28').Value = '9.5 MJ/m2'
$ws.Range('O28').Value = '10.5 °C'
$ws.Range('E29').Value = '2026-02-11 16:20:05'
$ws.Range('H29').Value = "'85%"
$ws.Range('K29').Value = '10.5 MJ/m2'
$ws.Range('O29').Value = '12.7 °C'
$ws.Range('E30').Value = '2026-02-11 16:20:08'
$ws.Range('K30').Value = '10.5 MJ/m2'
$ws.Range('O30').Value = '12.0 °C'
$ws.Range('E31').Value = '2026-02-11 16:20:11'
$ws.Range('K31').Value = '9.4 MJ/m2'
$ws.Range('E32').Value = '2026-02-11 16:20:13'
$ws.Range('H32').Value = "'71%"
$ws.Range('I32').Value = '2.6 mm'
$ws.Range('N32').Value = '8.4 °C 15:54 TU'
$ws.Range('O32').Value = '10.3 °C'
$ws.Range('E33').Value = '2026-02-11 16:20:16'
$ws.Range('H33').Value = "'81%"
$ws.Range('J33').Value = '1005.3 hPa'
$ws.Range('K33').Value = '10.7 MJ/m2'
$ws.Range('O33').Value = '6.6 °C'
$ws.Range('E34').Value = '2026-02-11 16:20:19'
$ws.Range('K34').Value = '10.9 MJ/m2'
$ws.Range('E35').Value = '2026-02-11 16:20:22'
$ws.Range('J35').Value = '1008.5 hPa'
$ws.Range('E36').Value = '2026-02-11 16:20:25'
$ws.Range('K36').Value = '10.5 MJ/m2'
$ws.Range('L36').Value = '29.5 km/h - 114º 15:33 TU'
$ws.Range('O36').Value = '12.8 °C'
$ws.Range('E37').Value = '2026-02-11 16:20:28'
$ws.Range('H37').Value = "'84%"
$ws.Range('J37').Value = '1005.1 hPa'
$ws.Range('L37').Value = '35.3 km/h - 229º 15:58 TU'
$ws.Range('O37').Value = '9.1 °C'
$ws.Range('E38').Value = '2026-02-11 16:20:30'
$ws.Range('H38').Value = "'60%"
$ws.Range('K38').Value = '10.7 MJ/m2'
$ws.Range('O38').Value = '15.5 °C'
$ws.Range('E39').Value = '2026-02-11 16:20:33'
$ws.Range('K39').Value = '12.4 MJ/m2'
$ws.Range('M39').Value = '4.7 °C 15:42 TU'
$ws.Range('O39').Value = '0.8 °C'
$ws.Range('E40').Value = '2026-02-11 16:20:36'
$ws.Range('H40').Value = "'89%"
$ws.Range('I40').Value = '0.3 mm'
$ws.Range('J40').Value = '1007.6 hPa'
$ws.Range('O40').Value = '7.4 °C'
$ws.Range('E41').Value = '2026-02-11 16:20:38'
$ws.Range('H41').Value = "'42%"
$ws.Range('K41').Value = '8.9 MJ/m2'
$ws.Range('O41').Value = '19.5 °C'
$ws.Range('E42').Value = '2026-02-11 16:20:41'
$ws.Range('H42').Value = "'90%"
$ws.Range('O42').Value = '12.4 °C'
$ws.Range('E43').Value = '2026-02-11 16:20:44'
$ws.Range('H43').Value = "'62%"
$ws.Range('K43').Value = '9.5 MJ/m2'
$ws.Range('O43').Value = '13.3 °C'
$ws.Range('E44').Value = '2026-02-11 16:20:46'
$ws.Range('H44').Value = "'81%"
$ws.Range('K44').Value = '7.9 MJ/m2'
$ws.Range('O44').Value = '-0.4 °C'
$ws.Range('E45').Value = '2026-02-11 16:20:49'
$ws.Range('G45').Value = '2 cm'
$ws.Range('H45').Value = "'87%"
$ws.Range('J45').Value = '1006.8 hPa'
$ws.Range('K45').Value = '6.2 MJ/m2'
$ws.Range('O45').Value = '6.6 °C'
$ws.Range('E46').Value = '2026-02-11 16:20:52'
$ws.Range('H46').Value = "'52%"
$ws.Range('J46').Value = '1008.1 hPa'
$ws.Range('K46').Value = '7.4 MJ/m2'
$ws.Range('N46').Value = '14.1 °C 15:54 TU'
$ws.Range('O46').Value = '17.9 °C'
